$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New shared strings must be created in the same order they first appear in
# the edited file (index 143 = "failed login attempt...", 144 = "app wouldn't
# run...", 145 = "Sonya/Ryan", 146 = "wrong header displays...", 147 = "if you
# log out..."), so set cells in that exact sequence.

# Add new remediation rows 67-70
$ws.Range("A67").Value = 16
$ws.Range("C67").Value = "failed login attempt routes to a blank page that only says try again or make an account"
$ws.Range("D67").Value = 1
$ws.Range("E67").Value = "Sarah"

$ws.Range("A68").Value = 17
$ws.Range("C68").Value = "app wouldn't run because of faulty html syntax"
$ws.Range("D68").Value = 1
$ws.Range("E68").Value = "Sonya"

# Update the existing "Assigned" entry on row 65 (E65) from "Sonya" to "Sonya/Ryan"
$ws.Range("E65").Value = "Sonya/Ryan"

$ws.Range("A69").Value = 18
$ws.Range("C69").Value = "wrong header displays on login, register, moodchoose, likes, dislikes pages"
$ws.Range("D69").Value = 1
$ws.Range("E69").Value = "Sonya"

$ws.Range("A70").Value = 19
$ws.Range("C70").Value = "if you log out, and then attempt to log back in but put the wrong password, internal service error"
$ws.Range("D70").Value = 2
$ws.Range("E70").Value = "Sarah"

# Update the view to match the new scroll/selection position
$ws.Range("C65").Select()
try {
    $excel.ActiveWindow.ScrollRow = 56
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scroll position is a cosmetic view-state setting; ignore if unsupported.
}
